# Updated capital structure database
# Applies the capital-structure data refresh for the Chile / Insurance (General)
# rows (row 2 = industry aggregate, row 3 = Banvida S.A.) per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: fix company name (remove comma) ----
$ws.Range("B3").Value = "Banvida S.A. (SNSE:BANVIDA)"

# ---- Row 2 value updates ----
$ws.Range("E2").Value = 0.149
$ws.Range("K2").Value = 84.59999999999999
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 0.02782608695652174
$ws.Range("W2").Value = 0.1579242113123016
$ws.Range("X2").Value = 0.04928836204419897
$ws.Range("Y2").Value = 0.1086358492681027
$ws.Range("AA2").Value = -0.001621754042136207
$ws.Range("AB2").Value = 0.04494648126111223
$ws.Range("AC2").Value = -0.04656823530324844
$ws.Range("AD2").Value = 94.59999999999999
$ws.Range("AF2").Value = 94.59999999999999
$ws.Range("AG2").Value = 78.59999999999999
$ws.Range("AH2").Value = 0.1412783751493429
$ws.Range("AI2").Value = 0.1461003861003861
$ws.Range("AJ2").Value = 0.1202570379436964
$ws.Range("AK2").Value = 0.1244655581947743
$ws.Range("AL2").Value = 4.64
$ws.Range("AM2").Value = 4.64
$ws.Range("AO2").Value = -0.2140086206896552
$ws.Range("AQ2").Value = -0.2140086206896552

# ---- Row 3 value updates ----
$ws.Range("E3").Value = 0.149
$ws.Range("K3").Value = 84.59999999999999
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 16
$ws.Range("V3").Value = 0.02782608695652174
$ws.Range("W3").Value = 0.1579242113123016
$ws.Range("X3").Value = 0.04928836204419897
$ws.Range("Y3").Value = 0.1086358492681027
$ws.Range("AA3").Value = -0.001621754042136207
$ws.Range("AB3").Value = 0.04494648126111223
$ws.Range("AC3").Value = -0.04656823530324844
$ws.Range("AD3").Value = 94.59999999999999
$ws.Range("AF3").Value = 94.59999999999999
$ws.Range("AG3").Value = 78.59999999999999
$ws.Range("AH3").Value = 0.1412783751493429
$ws.Range("AI3").Value = 0.1461003861003861
$ws.Range("AJ3").Value = 0.1202570379436964
$ws.Range("AK3").Value = 0.1244655581947743
$ws.Range("AL3").Value = 4.64
$ws.Range("AM3").Value = 4.64
$ws.Range("AO3").Value = -0.2140086206896552
$ws.Range("AQ3").Value = -0.2140086206896552
